$d = $word.ActiveDocument

# --- 1. Re-scope the "some-block-quotes-in-different-ways" bookmark so it
#        wraps the whole section (from before the heading run through to
#        the end of the final paragraph) instead of just the heading text.
$rng = $d.Range(0, $d.Content.End)
$d.Bookmarks.Add("some-block-quotes-in-different-ways", $rng)

# --- 2. Add the new "Section Number" character style (docx writer support
#        for --number-sections), based on the existing "Body Text Char"
#        style, matching the style used to number headings/sections.
$sectionNumber = $d.Styles.Add("Section Number", 2)
$sectionNumber.BaseStyle = "BodyTextChar"
